$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.291.93"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "2.189.84"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").Value = "'255.69"
$ws.Range("E5").Value = "  +3.83%  "

$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").Value = "'68.00"
$ws.Range("E7").Value = "  -3.99%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +4.64%  "

$ws.Range("D10").Value = "'37.94"
$ws.Range("E10").Value = "  +4.16%  "

$ws.Range("D11").Value = "'58.10"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "'0.0943"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").Value = "'7.09"
$ws.Range("E13").Value = "  +5.15%  "

$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "2.523.16"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "'0.870"
$ws.Range("E16").Value = "  +2.80%  "

$ws.Range("D17").Value = "'14.54"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "2.190.13"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").Value = "41.295.69"
$ws.Range("E19").Value = "  -0.95%  "

$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("E21").Value = "  +2.23%  "

$ws.Range("D22").Value = "'71.85"
$ws.Range("E22").Value = "  -2.24%  "

$ws.Range("D23").Value = "'232.52"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "'11.97"
$ws.Range("E25").Value = "  +20.51%  "

$ws.Range("D26").Value = "'3.84"
$ws.Range("E26").Value = "  +6.56%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'2.54"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").Value = "'169.96"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "'20.61"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'5.48"
$ws.Range("E34").Value = "  +6.29%  "

$ws.Range("D35").Value = "'0.0727"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("D36").Value = "'4.62"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").Value = "'25.27"
$ws.Range("E37").Value = "  +8.37%  "

$ws.Range("D38").Value = "'3.98"
$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("D39").Value = "'0.0298"
$ws.Range("E39").Value = "  +9.00%  "

$ws.Range("D40").Value = "'2.23"
$ws.Range("E40").Value = "  -2.35%  "

$ws.Range("D41").Value = "'5.76"
$ws.Range("E41").Value = "  -2.04%  "

$ws.Range("D42").Value = "'12.20"
$ws.Range("E42").Value = "  +18.15%  "

$ws.Range("D43").Value = "'64.29"
$ws.Range("E43").Value = "  -3.03%  "

$ws.Range("D44").Value = "'0.205"
$ws.Range("E44").Value = "  +7.70%  "

$ws.Range("D45").Value = "'4.89"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").Value = "'8.62"
$ws.Range("E47").Value = "  -3.92%  "

$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").Value = "'1.14"
$ws.Range("E49").Value = "  +3.69%  "

$ws.Range("E50").Value = "  -0.78%  "

$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  +0.38%  "
